# Update "想去人数" (want-to-go count, column F) figures across all four
# sheets to the newly scraped totals.
#
# Sheet "展览" (Exhibitions)
$ws = $excel.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 2364
$ws.Cells.Item(9, 6).Value = 2981
$ws.Cells.Item(10, 6).Value = 172
$ws.Cells.Item(11, 6).Value = 4442
$ws.Cells.Item(12, 6).Value = 389
$ws.Cells.Item(13, 6).Value = 214
$ws.Cells.Item(17, 6).Value = 138
$ws.Cells.Item(18, 6).Value = 221
$ws.Cells.Item(20, 6).Value = 109
$ws.Cells.Item(22, 6).Value = 4496
$ws.Cells.Item(24, 6).Value = 3775
$ws.Cells.Item(25, 6).Value = 1140
$ws.Cells.Item(26, 6).Value = 216
$ws.Cells.Item(27, 6).Value = 570
$ws.Cells.Item(30, 6).Value = 592
$ws.Cells.Item(31, 6).Value = 585
$ws.Cells.Item(32, 6).Value = 549

# Sheet "演出" (Performances)
$ws = $excel.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 32

# Sheet "本地生活" (Local life)
$ws = $excel.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 17

# Sheet "全部类型" (All types - aggregate)
$ws = $excel.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 17
$ws.Cells.Item(9, 6).Value = 2364
$ws.Cells.Item(13, 6).Value = 2981
$ws.Cells.Item(14, 6).Value = 172
$ws.Cells.Item(15, 6).Value = 4442
$ws.Cells.Item(16, 6).Value = 389
$ws.Cells.Item(17, 6).Value = 214
$ws.Cells.Item(21, 6).Value = 140
$ws.Cells.Item(22, 6).Value = 221
$ws.Cells.Item(25, 6).Value = 109
$ws.Cells.Item(27, 6).Value = 4496
$ws.Cells.Item(29, 6).Value = 3777
$ws.Cells.Item(30, 6).Value = 1140
$ws.Cells.Item(31, 6).Value = 216
$ws.Cells.Item(32, 6).Value = 570
$ws.Cells.Item(35, 6).Value = 592
$ws.Cells.Item(36, 6).Value = 585
$ws.Cells.Item(37, 6).Value = 549
$ws.Cells.Item(39, 6).Value = 32
